$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 611, shifting the existing
# rows 611:691 down to 613:693 (dimension grows from R691 to R693).
$ws.Rows("611:612").Insert()

# New row 611 (brand new price record)
$ws.Range("A611").Value = 3
$ws.Range("B611").Value = 'Femacal de La Calera'
$ws.Range("C611").Value = 'Coquimbo'
$ws.Range("D611").Value = 44984
$ws.Range("E611").Value = 5
$ws.Range("F611").Value = 100112003
$ws.Range("G611").Value = 'Ajo'
$ws.Range("H611").Value = 'Chino'
$ws.Range("I611").Value = 'Primera'
$ws.Range("J611").Value = 73
$ws.Range("K611").Value = 15000
$ws.Range("L611").Value = 15500
$ws.Range("M611").Value = 15260
$ws.Range("N611").Value = '$/caja 10 kilos'
$ws.Range("O611").Value = 'China'
$ws.Range("P611").Value = 1526
$ws.Range("Q611").Value = 10
$ws.Range("R611").Value = 'Hortaliza'

# New row 612 (brand new price record)
$ws.Range("A612").Value = 3
$ws.Range("B612").Value = 'Femacal de La Calera'
$ws.Range("C612").Value = 'Coquimbo'
$ws.Range("D612").Value = 44984
$ws.Range("E612").Value = 5
$ws.Range("F612").Value = 100112003
$ws.Range("G612").Value = 'Ajo'
$ws.Range("H612").Value = 'Chino'
$ws.Range("I612").Value = 'Primera'
$ws.Range("J612").Value = 75
$ws.Range("K612").Value = 18000
$ws.Range("L612").Value = 18500
$ws.Range("M612").Value = 18233
$ws.Range("N612").Value = '$/malla 10 kilos'
$ws.Range("O612").Value = 'China'
$ws.Range("P612").Value = 1823
$ws.Range("Q612").Value = 10
$ws.Range("R612").Value = 'Hortaliza'
